$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay as Text so values like "602.86" or "7.91"
# are not auto-converted into numbers by Excel (matches original inlineStr data).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.700.35"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "3.474.57"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "602.86"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("D6").Value = "147.29"
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("D7").Value = "3.472.46"
$ws.Range("E7").Value = "  -2.12%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  -2.96%  "
$ws.Range("D11").Value = "7.54"
$ws.Range("E11").Value = "  +3.07%  "
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").Value = "0.0000213"
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").Value = "4.069.32"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "31.36"
$ws.Range("E15").Value = "  -5.24%  "
$ws.Range("D16").Value = "3.469.97"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "66.715.68"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "6.42"
$ws.Range("E19").Value = "  -5.13%  "
$ws.Range("D20").Value = "15.33"
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("D21").Value = "10.02"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").Value = "438.08"
$ws.Range("E22").Value = "  -3.65%  "
$ws.Range("D23").Value = "0.607"
$ws.Range("E23").Value = "  -5.20%  "
$ws.Range("D24").Value = "79.39"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "3.619.41"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  -7.82%  "
$ws.Range("D28").Value = "9.79"
$ws.Range("E28").Value = "  -6.93%  "
$ws.Range("D29").Value = "8.37"
$ws.Range("E29").Value = "  -7.66%  "
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("E31").Value = "  -5.92%  "
$ws.Range("D32").Value = "0.167"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").Value = "25.33"
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").Value = "3.471.49"
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  -5.66%  "
$ws.Range("E37").Value = "  -6.22%  "
$ws.Range("B38").Value = "USDe"
$ws.Range("C38").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "7.91"
$ws.Range("E39").Value = "  -4.14%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "176.07"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "0.0882"
$ws.Range("E42").Value = "  -3.87%  "
$ws.Range("D43").Value = "2.11"
$ws.Range("E43").Value = "  -10.45%  "
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("D45").Value = "0.891"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "46.35"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").Value = "28.83"
$ws.Range("E47").Value = "  -6.35%  "
$ws.Range("E48").Value = "  -7.90%  "
$ws.Range("D49").Value = "7.44"
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("D50").Value = "2.42"
$ws.Range("E50").Value = "  -8.33%  "
$ws.Range("D51").Value = "0.978"
$ws.Range("E51").Value = "  -4.16%  "
